# Generate Report for Handoff
#
# Refresh the "Latest Handoff Datetime" column for every localization file
# that was just (re-)handed off, on each per-locale status sheet.
# Rows whose handoff datetime does not move (already reported handoff time
# is still current) are left untouched.

$wb = $excel.ActiveWorkbook

$handoffTimestamps = @{
    "zh-cn" = "2016-03-09 07:56:05"
    "de-de" = "2016-03-09 07:56:11"
}

# Rows on each locale sheet whose handoff just completed and therefore get
# a fresh "Latest Handoff Datetime" (column D) stamp.
$rowsToStamp = @(7, 10, 11, 12, 13, 14, 16)

foreach ($localeName in $handoffTimestamps.Keys) {
    $ws = $wb.Worksheets.Item($localeName)
    $timestamp = $handoffTimestamps[$localeName]

    foreach ($row in $rowsToStamp) {
        $ws.Cells.Item($row, 4).Value = $timestamp
    }
}
